$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A23").Value = "Java"
$ws.Range("B23").Value = "Java Command Line"
$ws.Range("C23").Value = "To call a java program on a command line, there are 2 ways:`n1. Call the runnable jar: >java -jar NameOfJar.jar arg1 arg2...`n2. Call the class: >java -cp {jar_path} com.myles.ClassName `narg1 arg2..."
$ws.Range("A23:C23").RowHeight = 33

